$d = $word.ActiveDocument

function FindParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t -like $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Split "Assertions where appropriate (for 'impossible' cases)." into
#    "... (for 'impossible' " + "or rare " + "cases)."
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Assertions where appropriate*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xml:space="preserve"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Assertions where appropriate</w:t></w:r><w:r><w:t xml:space="preserve"> (for ‘impossible’ </w:t></w:r><w:r><w:t xml:space="preserve">or rare </w:t></w:r><w:r><w:t>cases).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new sub-bullet right after the paragraph above:
#    "Especially check for things like integer overflows that would
#    otherwise be hard to debug." at ilvl=2 (ListLevelNumber=3). It takes
#    over the <w:lastRenderedPageBreak/> that used to sit on "Support move
#    semantics..." and the _GoBack bookmark that used to sit at the very
#    end of this list (on "Investigate a template header-only approach.").
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Improved error reporting via exceptions.*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.InsertParagraphBefore()
$newp = $d.Paragraphs.Item($idx)
$nr = $newp.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Especially check fo</w:t></w:r><w:r><w:t>r things like integer overflows that would otherwise be hard to debug.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$nr.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. Remove the <w:lastRenderedPageBreak/> that used to be on "Support move
#    semantics where appropriate." (it moved to the new bullet above).
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Support move semantics where appropriate.*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Support move semantics where appropriate.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4. Remove the _GoBack bookmark from the end of "Investigate a template
#    header-only approach." (it moved to the new bullet inserted in step 2).
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Investigate a template header-only approach.*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Investigate a template header-only approach.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5. "New Modules" heading gains a <w:lastRenderedPageBreak/>; the next
#    bullet ("Helper service to run HadesMem tools...") loses the one it had.
# ---------------------------------------------------------------------------
$idx = FindParaIndex("New Modules*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>New Modules</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

$idx = FindParaIndex("Helper service to run HadesMem*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xml:space="preserve"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Helper service to run HadesMem tools as ‘SYSTEM’ for when manipulation certain protected/critical processes (running in separate desktops, sessions, etc.).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6. "Patcher" heading gains a <w:lastRenderedPageBreak/>; the next bullet
#    ("VEH hooking (both INT3 and DR).") loses the one it had.
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Patcher*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Patcher</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

$idx = FindParaIndex("VEH hooking*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>VEH hooking (both INT3 and DR).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 7. "Load config directory." gains a <w:lastRenderedPageBreak/>; "Bound
#    import directory." loses the one it had.
# ---------------------------------------------------------------------------
$idx = FindParaIndex("Load config directory.*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Load config directory.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

$idx = FindParaIndex("Bound import directory.*")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Bound import directory.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Output "Done"
